$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42 / 43 swap (Filecoin <-> RenderToken) ---
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "5.34"
$ws.Range("E42").Value = "  +4.20%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.60"
$ws.Range("E43").Value = "  +4.99%  "

# --- Remaining cell value updates ---
$ws.Range("D2").Value = "56.748.73"
$ws.Range("E2").Value = "  +3.19%  "
$ws.Range("D3").Value = "2.325.22"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "520.67"
$ws.Range("E5").Value = "  +2.93%  "
$ws.Range("D6").Value = "134.57"
$ws.Range("E6").Value = "  +3.97%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "0.537"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("D9").Value = "2.352.04"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("E10").Value = "  +5.08%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "5.31"
$ws.Range("E12").Value = "  +4.02%  "
$ws.Range("D13").Value = "0.343"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").Value = "23.88"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "2.744.13"
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").Value = "56.822.01"
$ws.Range("E16").Value = "  +3.27%  "
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").Value = "2.324.07"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").Value = "10.48"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "4.22"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").Value = "322.72"
$ws.Range("E21").Value = "  +3.20%  "
$ws.Range("D22").Value = "6.56"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "60.83"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("E25").Value = "  +7.55%  "
$ws.Range("D26").Value = "0.994"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "7.99"
$ws.Range("E27").Value = "  +6.46%  "
$ws.Range("D28").Value = "1.29"
$ws.Range("E28").Value = "  +11.97%  "
$ws.Range("D29").Value = "0.0₃0740"
$ws.Range("E29").Value = "  +4.31%  "
$ws.Range("D30").Value = "1.71"
$ws.Range("E30").Value = "  +4.47%  "
$ws.Range("D31").Value = "167.09"
$ws.Range("E31").Value = "  -3.07%  "
$ws.Range("D32").Value = "6.19"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "18.34"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("D35").Value = "0.992"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").Value = "0.925"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "4.03"
$ws.Range("E38").Value = "  +3.41%  "
$ws.Range("E39").Value = "  +7.96%  "
$ws.Range("D40").Value = "37.97"
$ws.Range("E40").Value = "  +3.00%  "
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D44").Value = "138.68"
$ws.Range("E44").Value = "  +2.40%  "
$ws.Range("D45").Value = "281.95"
$ws.Range("E45").Value = "  +8.00%  "
$ws.Range("E46").Value = "  +2.44%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "0.565"
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("E49").Value = "  +2.72%  "
$ws.Range("D50").Value = "0.383"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("D51").Value = "17.78"
$ws.Range("E51").Value = "  +7.78%  "
